# Applies the cell-value updates for the "Updated cryptos list" commit.
# D-column price cells that look like plain numbers need NumberFormat forced
# to Text first so Excel keeps storing them as text (matching the source data)
# instead of silently converting them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.504.37"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "2.073.64"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.27"
$ws.Range("E5").Value = "  -1.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.628"
$ws.Range("E6").Value = "  +1.53%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.54"
$ws.Range("E8").Value = "  -1.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.392"
$ws.Range("E9").Value = "  +2.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0771"
$ws.Range("E10").Value = "  +0.90%  "
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("D12").Value = "2.379.04"
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.45"
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.69"
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.779"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.18"
$ws.Range("E16").Value = "  -0.88%  "
$ws.Range("D17").Value = "2.074.04"
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("D18").Value = "37.435.27"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.53"
$ws.Range("E19").Value = "  +4.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.73"
$ws.Range("E20").Value = "  +0.92%  "
$ws.Range("D21").Value = "0.0₃0817"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.60"
$ws.Range("E22").Value = "  +1.05%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.41"
$ws.Range("E25").Value = "  -1.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.18"
$ws.Range("E26").Value = "  +1.99%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.84"
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.45"
$ws.Range("E28").Value = "  +1.98%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.129"
$ws.Range("E29").Value = "  +1.04%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.17"
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.54"
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0623"
$ws.Range("E33").Value = "  -1.40%  "
$ws.Range("E34").Value = "  +2.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.50"
$ws.Range("E35").Value = "  -4.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.77"
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.28"
$ws.Range("E38").Value = "  -2.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.61"
$ws.Range("E39").Value = "  -4.71%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "99.23"
$ws.Range("E40").Value = "  +0.72%  "
$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.95"
$ws.Range("E41").Value = "  -0.65%  "
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0966"
$ws.Range("E42").Value = "  -1.98%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.482.65"
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.18"
$ws.Range("E44").Value = "  +3.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0213"
$ws.Range("E45").Value = "  +0.92%  "
$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.31"
$ws.Range("E46").Value = "  -5.86%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.03"
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.36"
$ws.Range("E48").Value = "  -4.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.23"
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.97"
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.74"
$ws.Range("E51").Value = "  +1.16%  "
